# Regenerate the localization-status report:
#  - "Ready for handoff" entries move to "In Translation"
#  - the affected status columns are re-fitted to the new (shorter) text

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$newWidth  = 13.4101845877511

# --- Overview sheet: zh-cn (col E) / de-de (col F) status columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

# --- zh-cn sheet: Status column (col C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C3").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth

# --- de-de sheet: Status column (col C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C3").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
